# "Create Product commit on 21/08/2023"
# Adds a new "DatProviderProd" worksheet (product lookup list) after the
# "Contact" sheet, tweaks a few Industry Type values on "DatProviderOrg" to
# reuse some of the newly introduced lookup values, and leaves
# "DatProviderOrg" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) DatProviderOrg: update a few Industry Type cells
# ---------------------------------------------------------------------
$orgSheet = $wb.Worksheets.Item("DatProviderOrg")
$orgSheet.Range("B2").Value = "Engineering"
$orgSheet.Range("B5").Value = "Consulting"
$orgSheet.Range("B6").Value = "Communications"

# ---------------------------------------------------------------------
# 2) Add the new "DatProviderProd" sheet, right after "Contact"
# ---------------------------------------------------------------------
$contactSheet = $wb.Worksheets.Item("Contact")
$prodSheet = $wb.Worksheets.Add($null, $contactSheet)
$prodSheet.Name = "DatProviderProd"

# Header cell: bold / filled / bordered, like the other sheets' headers,
# but left-aligned (copy the fill+border look from a non-bold header cell
# and just turn on bold so the alignment stays at "general").
$orgFirstSheet = $wb.Worksheets.Item("Organization")
$orgFirstSheet.Range("A1").Copy()
$prodSheet.Range("A1").PasteSpecial(-4122)
$prodSheet.Range("A1").Font.Bold = $true
$prodSheet.Range("A1").Value = "Product Name"

# Data rows: plain bordered style, matching the other lookup lists.
$orgSheet.Range("A2").Copy()
$prodSheet.Range("A2:A4").PasteSpecial(-4122)

$prodSheet.Range("A2").Value = "Hardware"
$prodSheet.Range("A3").Value = "Software"
$prodSheet.Range("A4").Value = "CRM Applications"

$prodSheet.Columns("A:A").ColumnWidth = 14.6

$prodSheet.Range("C2").Select()

# ---------------------------------------------------------------------
# 3) Leave "DatProviderOrg" as the active sheet/selection
# ---------------------------------------------------------------------
$orgSheet.Activate()
$orgSheet.Range("B2").Select()
